$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds free-form, numeric-looking text, e.g. "59.283.15"
# (thousands-grouped with dots), "519.50" (significant trailing zero), or
# "0.0₃0827" (subscript-compressed leading zeros). Setting NumberFormat to "@"
# (Text) on each target cell immediately before writing its value stops COM
# from coercing the string into a Double, which would otherwise normalize
# "519.50" -> 519.5, drop a trailing zero, or introduce float noise such as
# 6.2599999999999998.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.283.15"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.710.35"
$ws.Range("E3").Value = "  +6.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.50"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.22"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.706.45"
$ws.Range("E9").Value = "  +6.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.26"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E11").Value = "  +6.48%  "
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.143.82"
$ws.Range("E14").Value = "  +5.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.255.58"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.23"
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("E17").Value = "  +3.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.700.52"
$ws.Range("E18").Value = "  +6.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "356.20"
$ws.Range("E19").Value = "  +6.83%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.50"
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.23"
$ws.Range("E22").Value = "  +5.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.58"
$ws.Range("E24").Value = "  +3.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.424"
$ws.Range("E25").Value = "  +4.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.802.94"
$ws.Range("E26").Value = "  +5.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.162"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.990"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0827"
$ws.Range("E29").Value = "  +5.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.28"
$ws.Range("E30").Value = "  +5.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.996"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.46"
$ws.Range("E32").Value = "  +11.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.17"
$ws.Range("E33").Value = "  +3.31%  "
$ws.Range("E34").Value = "  +3.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "150.51"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("E36").Value = "  +11.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.08"
$ws.Range("E37").Value = "  +4.64%  "
$ws.Range("E38").Value = "  +4.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.862"
$ws.Range("E39").Value = "  +4.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.93"
$ws.Range("E40").Value = "  +3.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.72"
$ws.Range("E41").Value = "  +5.50%  "
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.624"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "282.44"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.08"
$ws.Range("E45").Value = "  +8.34%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0988"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0537"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0233"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.75"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.015.42"
$ws.Range("E51").Value = "  +6.75%  "
